$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldPrefix = "/host/d/GitHub/"
$newPrefix = "/host/c/Users/Junzhe/Desktop/Whole_heart_segmentation/GitHub/"

for ($row = 2; $row -le 12; $row++) {
    foreach ($col in @(3, 4)) {
        $cell = $ws.Cells.Item($row, $col)
        $text = [string]$cell.Value2
        $cell.Value2 = $text.Replace($oldPrefix, $newPrefix)
    }
}
